# Refresh the "cryptos" price list (sheet1) with the latest scraped values.
# Columns: A=index, B=Coin, C=Link, D=Price, E=Volume(1h)
# All of these cells are stored as text in the workbook, so numeric-looking
# values are prefixed with a leading apostrophe to force Excel to keep them
# as text instead of auto-converting them to floating point numbers (which
# would corrupt values like "0.07746" or "1.000").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.321.19"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").Value = "1.876.18"

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'0.7112"
$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("D6").Value = "'242.29"
$ws.Range("E6").Value = "  +0.75%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  +1.00%  "

$ws.Range("D9").Value = "'0.07746"
$ws.Range("E9").Value = "  +0.16%  "

$ws.Range("D10").Value = "'25.05"
$ws.Range("E10").Value = "  +0.15%  "

$ws.Range("D11").Value = "'0.08464"
$ws.Range("E11").Value = "  +2.53%  "

$ws.Range("D12").Value = "1.889.22"
$ws.Range("E12").Value = "  +1.62%  "

$ws.Range("D13").Value = "'5.213"
$ws.Range("E13").Value = "  -0.38%  "

$ws.Range("D14").Value = "'0.7111"
$ws.Range("E14").Value = "  -0.38%  "

$ws.Range("D15").Value = "'91.42"
$ws.Range("E15").Value = "  +1.33%  "

$ws.Range("D16").Value = "29.341.74"
$ws.Range("E16").Value = "  +0.54%  "

$ws.Range("E17").Value = "  +6.07%  "

$ws.Range("E18").Value = "  +2.37%  "

$ws.Range("D19").Value = "'242.49"
$ws.Range("E19").Value = "  -0.61%  "

# Rows 20 and 21 swapped their coin/link content (Avalanche <-> WrappedliquidstakedEther2.0)
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.132.01"
$ws.Range("E20").Value = "  +0.97%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'13.25"
$ws.Range("E21").Value = "  +0.58%  "

$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").Value = "'7.858"
$ws.Range("E23").Value = "  -1.13%  "

$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").Value = "'0.1608"
$ws.Range("E25").Value = "  +1.20%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").Value = "'9.023"
$ws.Range("E27").Value = "  +1.09%  "

$ws.Range("D28").Value = "'18.48"
$ws.Range("E28").Value = "  +1.04%  "

$ws.Range("D29").Value = "'1.515"
$ws.Range("E29").Value = "  +1.43%  "

$ws.Range("D30").Value = "'4.405"
$ws.Range("E30").Value = "  +0.43%  "

$ws.Range("D31").Value = "'4.329"
$ws.Range("E31").Value = "  +4.53%  "

$ws.Range("D32").Value = "'1.278"
$ws.Range("E32").Value = "  -2.75%  "

$ws.Range("D33").Value = "'0.05249"
$ws.Range("E33").Value = "  +1.07%  "

$ws.Range("D34").Value = "'1.933"
$ws.Range("E34").Value = "  +1.31%  "

$ws.Range("D35").Value = "'1.177"
$ws.Range("E35").Value = "  +0.30%  "

$ws.Range("D36").Value = "'0.7400"
$ws.Range("E36").Value = "  +1.70%  "

$ws.Range("D37").Value = "'2.685"
$ws.Range("E37").Value = "  +0.29%  "

$ws.Range("D38").Value = "'0.01867"
$ws.Range("E38").Value = "  +0.80%  "

$ws.Range("D39").Value = "'2.729"
$ws.Range("E39").Value = "  +1.62%  "

$ws.Range("D40").Value = "1.174.26"
$ws.Range("E40").Value = "  +1.34%  "

$ws.Range("D41").Value = "'6.387"
$ws.Range("E41").Value = "  +4.80%  "

$ws.Range("D42").Value = "'73.00"
$ws.Range("E42").Value = "  +0.74%  "

$ws.Range("D43").Value = "'0.8866"
$ws.Range("E43").Value = "  -2.04%  "

$ws.Range("D44").Value = "'106.28"
$ws.Range("E44").Value = "  +4.57%  "

$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  +0.07%  "

$ws.Range("D46").Value = "2.029.54"
$ws.Range("E46").Value = "  +0.92%  "

$ws.Range("D47").Value = "'1.813"
$ws.Range("E47").Value = "  +2.54%  "

$ws.Range("E48").Value = "  -0.47%  "

$ws.Range("E49").Value = "  +1.59%  "

$ws.Range("D50").Value = "'9.393"
$ws.Range("E50").Value = "  +0.83%  "

$ws.Range("D51").Value = "'0.4311"
$ws.Range("E51").Value = "  +1.17%  "
